# The "Version" row in the first table holds the version number string
# (originally "1.1.4", built from two runs "1.1." and "4" followed by the
# auto-managed "_GoBack" bookmark). The edit rewrites it to "5.5.5", built
# from a run "5" (before the bookmark, in the same position the old "1.1."
# run occupied) followed by the bookmark and three new runs ".5", "." and
# "5" after it.
#
# Because new runs created via InsertAfter/InsertBefore do not carry the
# <w:lang w:val="en-US"/> run-formatting that the original runs have (and
# this runtime does not expose a working LanguageID/Font language setter),
# we rebuild the whole paragraph via Range.InsertXML with the exact target
# WordprocessingML markup, keeping the paragraph/bookmark identifiers intact.

$d = $word.ActiveDocument

# Locate the table row labelled "Version" and take its second cell, which
# holds the version value we need to update.
$table = $d.Tables.Item(1)
$versionCell = $null
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $label = $table.Cell($i, 1).Range.Text
    if ($label.StartsWith("Version")) {
        $versionCell = $table.Cell($i, 2)
        break
    }
}

$targetRange = $versionCell.Range

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w:rsidR="00355EE4" w:rsidRDefault="00A7395D" w:rsidP="003F4029">
            <w:pPr>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>5</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>.5</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>.</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>5</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$targetRange.InsertXML($xml)
